# Final Test Case data
# - Drop the extra "R4..R8,R20" columns (E:J) and the trailing "R25,R26"
#   columns (now shifted to I:J), leaving Country,R1,R2,R3,R21,R22,R23,R24.
# - The remaining "Rxx" columns (E:H after the deletions) get a black font
#   color applied.
# - Selection is left on the old R21:R24 range (K1:N5) as it was right
#   before the deletion, matching the saved selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the about-to-be-orphaned selection before the columns shift away.
$ws.Range("K1:N5").Select()

# Remove columns R4, R5, R6, R7, R8, R20 (columns E through J).
$ws.Range("E1:J1").EntireColumn.Delete()

# Remove the trailing columns R25, R26, which have shifted left to I:J.
$ws.Range("I1:J1").EntireColumn.Delete()

# Apply a black font color to the remaining R21-R24 columns (now E:H).
$ws.Range("E1:H5").Font.Color = 0
